# 4.0.3 model and data
# Expand the "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv" rows on the
# "Boolean" sheet into per-vehicle-type CSV file rows, add six trailing blank
# rows, and refresh the workbook's view/active-tab state.

$wb = $excel.ActiveWorkbook
$boolean = $wb.Worksheets.Item("Boolean")
$integer = $wb.Worksheets.Item("Integer")
$about = $wb.Worksheets.Item("About")

# --- Boolean sheet: split the single "BVTQaZ.csv" row into six rows -------
$null = $boolean.Activate()

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv". Insert 5 more rows below
# it (so rows 17-22 are available) and fill them with the split file names.
$null = $boolean.Rows.Item(18).Resize(5).Insert()

$boolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$boolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$boolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$boolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$boolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$boolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the insert, "trans/BVTStL/BVTStL.csv", "trans/PVTStL/PVTStL.csv" and
# "trans/SRPbVT/SRPbVT.csv" now sit at rows 23-25, followed by
# "trans/VTQaZ/VTQaZ.csv" at row 26. Split that row the same way.
$null = $boolean.Rows.Item(27).Resize(5).Insert()

$boolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$boolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$boolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$boolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$boolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$boolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# "trans/VTStFES/VTStFES.csv" now sits at row 32, followed by six new blank
# rows (33-38) reserved for future entries.
$null = $boolean.Rows.Item(33).Resize(6).Insert()

$null = $boolean.Range("A32").Select()

# --- Integer sheet: no content change, just refresh the selection ---------
$null = $integer.Activate()
$null = $integer.Range("A13").Select()

# --- About sheet: becomes the active tab on reopen -------------------------
$null = $about.Activate()
$null = $about.Range("A1").Select()
